# Update countries & provincias Spain
#
# Refreshes the COVID-19 snapshot figures for a handful of countries and
# updates the "last updated" timestamp banner in A1. Bahrain's total-case
# count overtook Venezuela's, so those two rows (the sheet is kept sorted
# descending by column B, "Casos totales") swap labels: the row that used
# to be Venezuela now shows Barein's (updated) figures, and the row that
# used to be Barein now shows Venezuela's (unchanged) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 00:28"

# Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes
$data = @{
  4   = @("Estados Unidos",      6781572, 32283, 4055266, 2526316, 0, 990, 199990)
  8   = @("Peru",                 738020,  4160,  580753,  126340, 0, 115,  30927)
  29  = @("Canada",                138582,  572,  121567,    7827, 0,   9,   9188)
  45  = @("Guatemala",              82684,  512,   71983,    7717, 0,  12,   2984)
  47  = @("Japon",                  75958,  301,   67831,    6676, 0,   9,   1451)
  53  = @("Barein",                 61643,  678,   54831,    6599, 0,   0,    213)
  54  = @("Venezuela",              61569,    0,   49371,   11704, 0,   0,    494)
  58  = @("Nigeria",                56478,   90,   44430,   10960, 0,   5,   1088)
  96  = @("Guinea",                 10111,   50,    9444,     604, 0,   0,     63)
  108 = @("Luxemburgo",              7284,   40,    6555,     605, 0,   0,    124)
  132 = @("Gambia",                  3428,   23,    1737,    1586, 0,   2,    105)
  166 = @("Republica del Chad",      1087,    2,     944,      62, 0,   0,     81)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
  }
}
